$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 598.62067
$ws.Range("I15").Value = 598.62067
$ws.Range("K15").Value = 1795.86201
$ws.Range("M15").Value = -1626.86201
$ws.Range("H40").Value = 4698.625
$ws.Range("J40").Value = 4698.625
$ws.Range("L40").Value = 4698.625
$ws.Range("N40").Value = -5048.625
$ws.Range("H62").Value = 11147.333
$ws.Range("I62").Value = 9499
$ws.Range("J62").Value = 14444
$ws.Range("K62").Value = 9499
$ws.Range("L62").Value = 14444
$ws.Range("M62").Value = -8875
$ws.Range("N62").Value = -15692
$ws.Range("H65").Value = 11147.333
$ws.Range("I65").Value = 9499
$ws.Range("J65").Value = 14444
$ws.Range("K65").Value = 47495
$ws.Range("L65").Value = 72220
$ws.Range("M65").Value = -44375
$ws.Range("N65").Value = -78460
$ws.Range("H100").Value = 2331.56
$ws.Range("J100").Value = 972.6667
$ws.Range("L100").Value = 972.6667
$ws.Range("N100").Value = -2054.6667
$ws.Range("H112").Value = 2655.6875
$ws.Range("I112").Value = 1263.6666
$ws.Range("J112").Value = 2976.923
$ws.Range("K112").Value = 3790.9998
$ws.Range("L112").Value = 8930.769
$ws.Range("M112").Value = -2682.9998
$ws.Range("N112").Value = -11146.769
$ws.Range("H131").Value = 3666.6667
$ws.Range("H132").Value = 3163.7144
$ws.Range("I132").Value = 1817
$ws.Range("K132").Value = 5451
$ws.Range("M132").Value = -2921
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 349.25
$ws.Range("I11").Value = 298
$ws.Range("K11").Value = 298
$ws.Range("M11").Value = -154
$ws.Range("H30").Value = 7504.5
$ws.Range("I30").Value = 9
$ws.Range("J30").Value = 15000
$ws.Range("K30").Value = 9
$ws.Range("L30").Value = 15000
$ws.Range("M30").Value = 141
$ws.Range("N30").Value = -15300
$ws.Range("H45").Value = 1780.125
$ws.Range("J45").Value = 1998.5
$ws.Range("L45").Value = 1998.5
$ws.Range("N45").Value = -2752.5
$ws.Range("H101").Value = 29999
$ws.Range("J101").Value = 29999
$ws.Range("L101").Value = 29999
$ws.Range("N101").Value = -36489
$ws.Range("H110").Value = 1445.7333
$ws.Range("I110").Value = 1206.6923
$ws.Range("J110").Value = 2999.5
$ws.Range("K110").Value = 1206.6923
$ws.Range("L110").Value = 2999.5
$ws.Range("M110").Value = 838.3077000000001
$ws.Range("N110").Value = -7089.5
$ws.Range("H122").Value = 1411.091
$ws.Range("I122").Value = 1411.091
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4233.272999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1783.272999999999
$ws.Range("N122").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 816
$ws.Range("I25").Value = 554.6667
$ws.Range("J25").Value = 1600
$ws.Range("K25").Value = 554.6667
$ws.Range("L25").Value = 1600
$ws.Range("M25").Value = -319.6667
$ws.Range("N25").Value = -2070
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H58").Value = 561.4
$ws.Range("I58").Value = 552.25
$ws.Range("J58").Value = 598
$ws.Range("K58").Value = 552.25
$ws.Range("L58").Value = 598
$ws.Range("M58").Value = -349.25
$ws.Range("N58").Value = -1004
$ws.Range("H88").Value = 23140.166
$ws.Range("J88").Value = 23140.166
$ws.Range("L88").Value = 23140.166
$ws.Range("N88").Value = -23952.166
$ws.Range("H91").Value = 23140.166
$ws.Range("J91").Value = 23140.166
$ws.Range("L91").Value = 23140.166
$ws.Range("N91").Value = -25948.166
$ws.Range("H132").Value = 7909.5713
$ws.Range("I132").Value = 2186.8
$ws.Range("K132").Value = 6560.400000000001
$ws.Range("M132").Value = -4030.400000000001
$ws.Range("H136").Value = 561.4
$ws.Range("I136").Value = 552.25
$ws.Range("J136").Value = 598
$ws.Range("K136").Value = 1656.75
$ws.Range("L136").Value = 1794
$ws.Range("M136").Value = 893.25
$ws.Range("N136").Value = -6894
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 177092.17
$ws.Range("J12").Value = 1137.3334
$ws.Range("L12").Value = 3412.0002
$ws.Range("N12").Value = -3758.0002
$ws.Range("H18").Value = 987.5
$ws.Range("I18").Value = 987.5
$ws.Range("K18").Value = 2962.5
$ws.Range("M18").Value = -2793.5
$ws.Range("H98").Value = 583.75
$ws.Range("I98").Value = 546.75
$ws.Range("J98").Value = 620.75
$ws.Range("K98").Value = 1640.25
$ws.Range("L98").Value = 1862.25
$ws.Range("M98").Value = -142.25
$ws.Range("N98").Value = -4858.25
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2899.5
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H122").Value = 334666.66
$ws.Range("I122").Value = 334666.66
$ws.Range("K122").Value = 1003999.98
$ws.Range("M122").Value = -1001549.98
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2490.2
$ws.Range("J68").Value = 3001.5
$ws.Range("L68").Value = 3001.5
$ws.Range("N68").Value = -4499.5
$ws.Range("H71").Value = 2490.2
$ws.Range("J71").Value = 3001.5
$ws.Range("L71").Value = 15007.5
$ws.Range("N71").Value = -22495.5
$ws.Range("H132").Value = 6080.6665
$ws.Range("I132").Value = 5533
$ws.Range("K132").Value = 16599
$ws.Range("M132").Value = -14069
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 29500
$ws.Range("J119").Value = 29000
$ws.Range("L119").Value = 29000
$ws.Range("N119").Value = -38676
$ws.Range("H122").Value = 4004
$ws.Range("I122").Value = 4004
$ws.Range("K122").Value = 12012
$ws.Range("M122").Value = -9562
$ws.Range("H126").Value = 4500
$ws.Range("I126").Value = 4000
$ws.Range("K126").Value = 12000
$ws.Range("M126").Value = -9530

Write-Host "Applied all edits"